$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2:53 (columns A:E) after inserting a new leading
# observation row and rewriting the forecast columns with the bugfixed
# (simulated) AR2 evaluation values.
$data = @(
    @(39400, 2007, -0.7015558851707571, 2008, $null),
    @(39583, 2008, $null, 2009, $null),
    @(39765, 2008, 0.3590181115727509, 2009, $null),
    @(39948, 2009, $null, 2010, $null),
    @(40130, 2009, -0.01587181126743165, 2010, $null),
    @(40310, 2010, 0.1037184923425638, 2011, 0.1523866428317699),
    @(40494, 2010, -0.02256889165885845, 2011, 0.1850158025574977),
    @(40676, 2011, -0.3429736775751868, 2012, -0.1149259972201833),
    @(40862, 2011, 0.09611428386597787, 2012, -0.0461580488825053),
    @(41044, 2012, -0.2064822032187164, 2013, -0.1075365560329056),
    @(41228, 2012, -0.1827723404408288, 2013, -0.1118170657869455),
    @(41409, 2013, 0.005797897989445744, 2014, -0.05125851421730054),
    @(41592, 2013, -0.001350220946483294, 2014, -0.02003709364999384),
    @(41774, 2014, 0.2273278155300318, 2015, 0.001779459905826286),
    @(41957, 2014, -0.075754880139145, 2015, -0.2383077634182995),
    @(42137, 2015, -0.2723669344146917, 2016, -0.1573879283727764),
    @(42321, 2015, -0.5761528471665445, 2016, -0.2318455351884685),
    @(42503, 2016, -0.05011389829933099, 2017, -0.1100689213476058),
    @(42689, 2016, -0.2011999787958185, 2017, -0.2555440101933648),
    @(42867, 2017, 0.2788833036191596, 2018, -0.1416348838281123),
    @(43053, 2017, 0.1213692818849532, 2018, -0.06802695342326137),
    @(43145, 2018, 0.1319881130530609, 2019, -0.1311012291692171),
    @(43235, 2018, -0.1459594536071473, 2019, -0.1825768856163368),
    @(43326, 2018, -0.003510046413479628, 2019, -0.05140618367761585),
    @(43418, 2018, 0.1493219406571766, 2019, 0.2006993856294326),
    @(43510, 2019, 0.3205506405107528, 2020, 0.02707231012755162),
    @(43600, 2019, -0.2305835819295887, 2020, -0.2087674606261247),
    @(43691, 2019, -0.3355560014489445, 2020, -0.3121234745279589),
    @(43783, 2019, -0.4278219446121612, 2020, -0.5866687821558636),
    @(43875, 2020, -0.4642429512678015, 2021, -0.18419261498126),
    @(43966, 2020, -0.3292346503903532, 2021, -0.09082057608673644),
    @(44068, 2020, -0.4412356890029168, 2021, -0.1505745581714235),
    @(44159, 2020, -1.026566979837418, 2021, -2.497698913894009),
    @(44251, 2021, -0.4955027878269735, 2022, -0.278034778427505),
    @(44341, 2021, -0.05652554630790618, 2022, 0.149910086121019),
    @(44432, 2021, -0.04976849661378902, 2022, -0.4624617512119622),
    @(44525, 2021, 0.3179894933462268, 2022, 0.1527842396192636),
    @(44617, 2022, 0.4286822041105198, 2023, -0.07547397514037213),
    @(44706, 2022, 1.183936177705625, 2023, 0.1147913851119675),
    @(44798, 2022, 0.4784173072842179, 2023, -0.05252235973735031),
    @(44890, 2022, 0.463604920919658, 2023, -0.2029074457040814),
    @(44981, 2023, 0.4661842108264924, 2024, -0.09280859516105355),
    @(45071, 2023, 0.8096199637471102, 2024, -0.08981539700775309),
    @(45163, 2023, 0.4279153732809959, 2024, -0.1676828594555424),
    @(45254, 2023, 0.621639092134818, 2024, -0.1003146131857147),
    @(45345, 2024, -0.438218070630314, 2025, -0.08844532960863294),
    @(45436, 2024, -1.085560181261136, 2025, -0.09271976299123352),
    @(45534, 2024, -0.9378224616154895, 2025, 0.0525099363110515),
    @(45618, 2024, -0.6768900623516982, 2025, 0.6778422458094902),
    @(45713, 2025, 1.664803039325657, 2026, 0.06351685110044247),
    @(45800, 2025, 0.5641976336596244, 2026, -0.1178589452312528),
    @(45891, 2025, 1.136769786738334, 2026, 0.06731424764072624)
)

# Clear the previous data block (old rows 2:52) before rewriting rows 2:53.
$ws.Range("A2:E52").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    if ($null -ne $row[0]) { $ws.Cells.Item($r, 1).Value2 = $row[0] }
    if ($null -ne $row[1]) { $ws.Cells.Item($r, 2).Value2 = $row[1] }
    if ($null -ne $row[2]) { $ws.Cells.Item($r, 3).Value2 = $row[2] }
    if ($null -ne $row[3]) { $ws.Cells.Item($r, 4).Value2 = $row[3] }
    if ($null -ne $row[4]) { $ws.Cells.Item($r, 5).Value2 = $row[4] }
}

# Row 53 is newly used by column A (date style "s=2"); copy that date
# number-format/font/border/alignment style from the row above (A52)
# without disturbing the value we already wrote into A53.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = 0
